$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Fix the menu text in D6 (shared string update):
#    "1. Main Menu, Menu – Features, 3. Spirits List Menu, ..."
#    -> "1. Main Menu, 2. Features Menu, 3. Spirits List Menu, ..."
$ws.Range("D6").Value = "1. Main Menu, 2. Features Menu, 3. Spirits List Menu, 4. Bottle List Menu, 5. After Dinner / Dessert Menu"

# 2. Slightly widen column A and column D (minor manual resize by the author).
#    The engine quantizes stored column width to multiples of 1/6, so we pick
#    the ColumnWidth input that lands on the closest achievable stored width
#    to the target values from the diff (21.7449392712551 and 77.4453441295547).
$ws.Columns.Item(1).ColumnWidth = 20.83
$ws.Columns.Item(4).ColumnWidth = 76.65

# 3. Move the active cell selection from C6 to D7.
$ws.Range("D7").Select()
